$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.302.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.367.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.87%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.40"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.95"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.50"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.85%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.729.07"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.41"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.352.79"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.244.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +14.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.48%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.60"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.14"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.85"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.53"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.03%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.31"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.43"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0950"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.26"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.78"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.87"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.69%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.73"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.98%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.74%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.09"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.42%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.866.47"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +13.07%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.227"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.60%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.74"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.17"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.19"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.23%  "
